$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update vm_pu results for the 380 kV case (B column target V set to 1.02 pu,
# columns C-F and I-N recalculated bus voltage magnitudes per row).

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.033575495278463
$ws.Range("D2").Value = 1.043705930972616
$ws.Range("E2").Value = 1.043365700692434
$ws.Range("F2").Value = 1.055482981067455
$ws.Range("I2").Value = 1.038939204012973
$ws.Range("J2").Value = 1.038699467832343
$ws.Range("K2").Value = 1.046479042136435
$ws.Range("L2").Value = 1.046139770677589
$ws.Range("M2").Value = 1.058223315619796
$ws.Range("N2").Value = 1.01687301748773

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.034441840983624
$ws.Range("D3").Value = 1.044384326505453
$ws.Range("E3").Value = 1.04412026524071
$ws.Range("F3").Value = 1.056286249143655
$ws.Range("I3").Value = 1.039108097939357
$ws.Range("J3").Value = 1.039209006071903
$ws.Range("K3").Value = 1.04696887389699
$ws.Range("L3").Value = 1.046705503148035
$ws.Range("M3").Value = 1.058840059859621
$ws.Range("N3").Value = 1.017042306086185

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.0350029481864
$ws.Range("D4").Value = 1.044823394183295
$ws.Range("E4").Value = 1.044609280838899
$ws.Range("F4").Value = 1.056806568387231
$ws.Range("I4").Value = 1.039215789032922
$ws.Range("J4").Value = 1.039538586840911
$ws.Range("K4").Value = 1.047285248314255
$ws.Range("L4").Value = 1.047071666756367
$ws.Range("M4").Value = 1.059239029646311
$ws.Range("N4").Value = 1.017151780011684

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.035238961474635
$ws.Range("D5").Value = 1.045008000436116
$ws.Range("E5").Value = 1.04481504359045
$ws.Range("F5").Value = 1.057025440583717
$ws.Range("I5").Value = 1.039260679738921
$ws.Range("J5").Value = 1.039677111509133
$ws.Range("K5").Value = 1.047418111954618
$ws.Range("L5").Value = 1.047225623704436
$ws.Range("M5").Value = 1.059406729780138
$ws.Range("N5").Value = 1.017197786326697

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.03527859636062
$ws.Range("D6").Value = 1.045038997886656
$ws.Range("E6").Value = 1.044849602599507
$ws.Range("F6").Value = 1.057062197754437
$ws.Range("I6").Value = 1.039268194622571
$ws.Range("J6").Value = 1.039700368550501
$ws.Range("K6").Value = 1.04744041208335
$ws.Range("L6").Value = 1.047251474988434
$ws.Range("M6").Value = 1.059434885753485
$ws.Range("N6").Value = 1.01720551000724

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.035006101322898
$ws.Range("D7").Value = 1.044825860816162
$ws.Range("E7").Value = 1.044612029543514
$ws.Range("F7").Value = 1.056809492459926
$ws.Range("I7").Value = 1.039216390369322
$ws.Range("J7").Value = 1.039540437938059
$ws.Range("K7").Value = 1.04728702419728
$ws.Range("L7").Value = 1.047073723852906
$ws.Range("M7").Value = 1.059241270570932
$ws.Range("N7").Value = 1.017152394815815

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.033868171530305
$ws.Range("D8").Value = 1.043935176027402
$ws.Range("E8").Value = 1.043620550480458
$ws.Range("F8").Value = 1.055754334001188
$ws.Range("I8").Value = 1.038996612140474
$ws.Range("J8").Value = 1.038871693690315
$ws.Range("K8").Value = 1.046644702198482
$ws.Range("L8").Value = 1.046330941701963
$ws.Range("M8").Value = 1.058431767979011
$ws.Range("N8").Value = 1.016930242986027

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.031867071594625
$ws.Range("D9").Value = 1.042366525522974
$ws.Range("E9").Value = 1.041879352237226
$ws.Range("F9").Value = 1.053899316379704
$ws.Range("I9").Value = 1.038597161323086
$ws.Range("J9").Value = 1.037692383061638
$ws.Range("K9").Value = 1.045508473792348
$ws.Range("L9").Value = 1.045022869343628
$ws.Range("M9").Value = 1.057004588954001
$ws.Range("N9").Value = 1.0165382891019

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.030535837085413
$ws.Range("D10").Value = 1.041321438669748
$ws.Range("E10").Value = 1.040722636899544
$ws.Range("F10").Value = 1.052665648585174
$ws.Range("I10").Value = 1.038322728244628
$ws.Range("J10").Value = 1.036905647528477
$ws.Range("K10").Value = 1.044748131494854
$ws.Range("L10").Value = 1.044151446475254
$ws.Range("M10").Value = 1.056052740746718
$ws.Range("N10").Value = 1.016276681337881

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.029960088454666
$ws.Range("D11").Value = 1.040869088677685
$ws.Range("E11").Value = 1.040222757418995
$ws.Range("F11").Value = 1.052132195095098
$ws.Range("I11").Value = 1.038201977655931
$ws.Range("J11").Value = 1.036564872766003
$ws.Range("K11").Value = 1.044418235072174
$ws.Range("L11").Value = 1.043774276779374
$ws.Range("M11").Value = 1.055640506422541
$ws.Range("N11").Value = 1.016163335555448

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.029746334142614
$ws.Range("D12").Value = 1.040701094550502
$ws.Range("E12").Value = 1.040037229855151
$ws.Range("F12").Value = 1.051934158758182
$ws.Range("I12").Value = 1.038156838011599
$ws.Range("J12").Value = 1.036438278111318
$ws.Range("K12").Value = 1.04429559862627
$ws.Range("L12").Value = 1.043634204939841
$ws.Range("M12").Value = 1.055487374175847
$ws.Range("N12").Value = 1.016121224135772

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.029792180421397
$ws.Range("D13").Value = 1.040737128528416
$ws.Range("E13").Value = 1.040077019333564
$ws.Range("F13").Value = 1.051976633123408
$ws.Range("I13").Value = 1.03816653361491
$ws.Range("J13").Value = 1.0364654337864
$ws.Range("K13").Value = 1.044321908982531
$ws.Range("L13").Value = 1.043664249635637
$ws.Range("M13").Value = 1.055520221992031
$ws.Range("N13").Value = 1.016130257612192

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.029942417318884
$ws.Range("D14").Value = 1.04085520163263
$ws.Range("E14").Value = 1.040207418583648
$ws.Range("F14").Value = 1.052115823052868
$ws.Range("I14").Value = 1.038198252258856
$ws.Range("J14").Value = 1.036554408726095
$ws.Range("K14").Value = 1.044408099900931
$ws.Range("L14").Value = 1.043762697863356
$ws.Range("M14").Value = 1.055627848659207
$ws.Range("N14").Value = 1.016159854808641

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.030034997103981
$ws.Range("D15").Value = 1.040927954240093
$ws.Range("E15").Value = 1.040287781786991
$ws.Range("F15").Value = 1.05220159746693
$ws.Range("I15").Value = 1.038217757091094
$ws.Range("J15").Value = 1.036609227075449
$ws.Range("K15").Value = 1.044461191996671
$ws.Range("L15").Value = 1.043823358523946
$ws.Range("M15").Value = 1.055694159715796
$ws.Range("N15").Value = 1.01617808933979

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.030574062179358
$ws.Range("D16").Value = 1.041351463587483
$ws.Range("E16").Value = 1.040755833180855
$ws.Range("F16").Value = 1.052701067741975
$ws.Range("I16").Value = 1.038330701699381
$ws.Range("J16").Value = 1.036928261351744
$ws.Range("K16").Value = 1.044770011781809
$ws.Range("L16").Value = 1.044176481523254
$ws.Range("M16").Value = 1.056080097881017
$ws.Range("N16").Value = 1.016284202332231

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.030912387828276
$ws.Range("D17").Value = 1.041617169167022
$ws.Range("E17").Value = 1.041049694687434
$ws.Range("F17").Value = 1.053014569817539
$ws.Range("I17").Value = 1.038401035627931
$ws.Range("J17").Value = 1.037128353735063
$ws.Range("K17").Value = 1.044963549716549
$ws.Range("L17").Value = 1.044398030666461
$ws.Range("M17").Value = 1.056322167019607
$ws.Range("N17").Value = 1.016350746350406

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.031109793418153
$ws.Range("D18").Value = 1.041772167933998
$ws.Range("E18").Value = 1.041221194135593
$ws.Range("F18").Value = 1.053197500849117
$ws.Range("I18").Value = 1.038441874981573
$ws.Range("J18").Value = 1.037245053241888
$ws.Range("K18").Value = 1.045076373036824
$ws.Range("L18").Value = 1.044527272084308
$ws.Range("M18").Value = 1.056463354240652
$ws.Range("N18").Value = 1.016389553769374

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.031177114765809
$ws.Range("D19").Value = 1.041825021369264
$ws.Range("E19").Value = 1.041279687074429
$ws.Range("F19").Value = 1.053259887514999
$ws.Range("I19").Value = 1.03845576870058
$ws.Range("J19").Value = 1.037284842854763
$ws.Range("K19").Value = 1.045114831960907
$ws.Range("L19").Value = 1.044571342662453
$ws.Range("M19").Value = 1.05651149409395
$ws.Range("N19").Value = 1.016402784966382

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.030876081854355
$ws.Range("D20").Value = 1.041588659680862
$ws.Range("E20").Value = 1.041018156289139
$ws.Range("F20").Value = 1.052980926697827
$ws.Range("I20").Value = 1.038393508616826
$ws.Range("J20").Value = 1.037106886844187
$ws.Range("K20").Value = 1.044942791545941
$ws.Range("L20").Value = 1.044374258921368
$ws.Range("M20").Value = 1.05629619606759
$ws.Range("N20").Value = 1.016343607484465

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.029898173404132
$ws.Range("D21").Value = 1.040820431236613
$ws.Range("E21").Value = 1.040169015101207
$ws.Range("F21").Value = 1.052074831979369
$ws.Range("I21").Value = 1.038188919838376
$ws.Range("J21").Value = 1.036528208248291
$ws.Range("K21").Value = 1.044382721523917
$ws.Range("L21").Value = 1.043733706591061
$ws.Range("M21").Value = 1.055596155553072
$ws.Range("N21").Value = 1.016151139437491

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.029283928111511
$ws.Range("D22").Value = 1.040337583057467
$ws.Range("E22").Value = 1.039635994690859
$ws.Range("F22").Value = 1.051505783463276
$ws.Range("I22").Value = 1.038058623788829
$ws.Range("J22").Value = 1.036164279629479
$ws.Range("K22").Value = 1.044030015662946
$ws.Range("L22").Value = 1.043331115765947
$ws.Range("M22").Value = 1.055155954869089
$ws.Range("N22").Value = 1.016030070953026

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.029609493357416
$ws.Range("D23").Value = 1.040593533522712
$ws.Range("E23").Value = 1.039918475981007
$ws.Range("F23").Value = 1.0518073846637
$ws.Range("I23").Value = 1.038127853537545
$ws.Range("J23").Value = 1.03635721322827
$ws.Range("K23").Value = 1.044217045045921
$ws.Range("L23").Value = 1.043544522103454
$ws.Range("M23").Value = 1.055389318484182
$ws.Range("N23").Value = 1.016094256822639

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.03089248675754
$ws.Range("D24").Value = 1.041601541837559
$ws.Range("E24").Value = 1.041032406840449
$ws.Range("F24").Value = 1.052996128357638
$ws.Range("I24").Value = 1.038396910321587
$ws.Range("J24").Value = 1.037116586842257
$ws.Range("K24").Value = 1.044952171467876
$ws.Range("L24").Value = 1.044385000301847
$ws.Range("M24").Value = 1.056307931246605
$ws.Range("N24").Value = 1.01634683325108

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.032383911339775
$ws.Range("D25").Value = 1.042771947336801
$ws.Range("E25").Value = 1.042328781786001
$ws.Range("F25").Value = 1.054378360676633
$ws.Range("I25").Value = 1.038701866163696
$ws.Range("J25").Value = 1.037997361765189
$ws.Range("K25").Value = 1.045802725774883
$ws.Range("L25").Value = 1.045360933339335
$ws.Range("M25").Value = 1.05737362571549
$ws.Range("N25").Value = 1.016639674334007

Write-Output "Updated vm_pu values for 380 kV case"